$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 21.50657928005314
$ws.Range("C2").Value = 4.807414580812353
$ws.Range("D2").Value = 13.75729255573399
$ws.Range("E2").Value = 13.95668731950672
$ws.Range("G2").Value = 65.79123157865165
$ws.Range("H2").Value = 23.54513563853234
$ws.Range("J2").Value = 8.419269586400425
$ws.Range("K2").Value = 16.9316321518543
$ws.Range("L2").Value = 12.97229514029918
$ws.Range("N2").Value = 25.01581923832431

$ws.Range("B3").Value = 21.39576014710026
$ws.Range("C3").Value = 4.691420675132188
$ws.Range("D3").Value = 13.75464716235898
$ws.Range("E3").Value = 13.97717114422052
$ws.Range("G3").Value = 65.7531867498833
$ws.Range("H3").Value = 23.57791439706084
$ws.Range("J3").Value = 8.428455755299744
$ws.Range("K3").Value = 16.85892844362725
$ws.Range("L3").Value = 12.98376566720636
$ws.Range("N3").Value = 25.0680551662164

$ws.Range("B4").Value = 21.33248486776053
$ws.Range("C4").Value = 4.617721529949431
$ws.Range("D4").Value = 13.75544402981763
$ws.Range("E4").Value = 13.9914972851779
$ws.Range("G4").Value = 65.74272764591976
$ws.Range("H4").Value = 23.60138614943526
$ws.Range("J4").Value = 8.434410679509906
$ws.Range("K4").Value = 16.81797042859998
$ws.Range("L4").Value = 12.99271717315033
$ws.Range("N4").Value = 25.10205994411956

$ws.Range("B5").Value = 21.30792007885186
$ws.Range("C5").Value = 4.587085432176855
$ws.Range("D5").Value = 13.75637861274496
$ws.Range("E5").Value = 13.99777540605372
$ws.Range("G5").Value = 65.74171131898726
$ws.Range("H5").Value = 23.61179182100224
$ws.Range("J5").Value = 8.436916716071757
$ws.Range("K5").Value = 16.8022188183211
$ws.Range("L5").Value = 12.99684527565888
$ws.Range("N5").Value = 25.11640348161173

$ws.Range("B6").Value = 21.30391541433221
$ws.Range("C6").Value = 4.581962522076659
$ws.Range("D6").Value = 13.75657065809027
$ws.Range("E6").Value = 13.9988444725439
$ws.Range("G6").Value = 65.74173861822989
$ws.Range("H6").Value = 23.61357043683952
$ws.Range("J6").Value = 8.43733764192295
$ws.Range("K6").Value = 16.79966035597322
$ws.Range("L6").Value = 12.99755976124853
$ws.Range("N6").Value = 25.11881460428313

$ws.Range("B7").Value = 21.33214861013067
$ws.Range("C7").Value = 4.617310775084998
$ws.Range("D7").Value = 13.75545416348768
$ws.Range("E7").Value = 13.99158017177649
$ws.Range("G7").Value = 65.74270079611188
$ws.Range("H7").Value = 23.6015230806214
$ws.Range("J7").Value = 8.434444155146933
$ws.Range("K7").Value = 16.81775417799343
$ws.Range("L7").Value = 12.99277090108921
$ws.Range("N7").Value = 25.102251416288

$ws.Range("B8").Value = 21.46739260611444
$ws.Range("C8").Value = 4.767945897692545
$ws.Range("D8").Value = 13.75587861192405
$ws.Range("E8").Value = 13.96338733137593
$ws.Range("G8").Value = 65.77543786913768
$ws.Range("H8").Value = 23.55574312945846
$ws.Range("J8").Value = 8.422371838931801
$ws.Range("K8").Value = 16.90580789694627
$ws.Range("L8").Value = 12.9758543044733
$ws.Range("N8").Value = 25.03342960364078

$ws.Range("B9").Value = 21.76941638232121
$ws.Range("C9").Value = 5.04286722200214
$ws.Range("D9").Value = 13.77586652844018
$ws.Range("E9").Value = 13.92196483305288
$ws.Range("G9").Value = 65.94184263097341
$ws.Range("H9").Value = 23.49253388749074
$ws.Range("J9").Value = 8.401182617566413
$ws.Range("K9").Value = 17.10706609295741
$ws.Range("L9").Value = 12.95780572296325
$ws.Range("N9").Value = 24.91377131072759

$ws.Range("B10").Value = 22.01227583409996
$ws.Range("C10").Value = 5.23139212290034
$ws.Range("D10").Value = 13.80214111240441
$ws.Range("E10").Value = 13.89996388817684
$ws.Range("G10").Value = 66.12610089388812
$ws.Range("H10").Value = 23.46231354194726
$ws.Range("J10").Value = 8.387113539256019
$ws.Range("K10").Value = 17.27143411428361
$ws.Range("L10").Value = 12.95373698822018
$ws.Range("N10").Value = 24.83514958652532

$ws.Range("B11").Value = 22.12696572881102
$ws.Range("C11").Value = 5.314055417169343
$ws.Range("D11").Value = 13.81658607318982
$ws.Range("E11").Value = 13.8917811912249
$ws.Range("G11").Value = 66.22328232606772
$ws.Range("H11").Value = 23.45209053885017
$ws.Range("J11").Value = 8.381035166462162
$ws.Range("K11").Value = 17.3495747245437
$ws.Range("L11").Value = 12.95387367001959
$ws.Range("N11").Value = 24.80139249511497

$ws.Range("B12").Value = 22.17097054974226
$ws.Range("C12").Value = 5.344899163478864
$ws.Range("D12").Value = 13.82241181939559
$ws.Range("E12").Value = 13.88894464944225
$ws.Range("G12").Value = 66.26199128626953
$ws.Range("H12").Value = 23.44872614922999
$ws.Range("J12").Value = 8.378779448396502
$ws.Range("K12").Value = 17.37962913750922
$ws.Range("L12").Value = 12.95421032596751
$ws.Range("N12").Value = 24.78889785672797

$ws.Range("B13").Value = 22.16146829550527
$ws.Range("C13").Value = 5.338277030426429
$ws.Range("D13").Value = 13.82114136409126
$ws.Range("E13").Value = 13.88954390170132
$ws.Range("G13").Value = 66.25356996616917
$ws.Range("H13").Value = 23.44942818970573
$ws.Range("J13").Value = 8.379263214151873
$ws.Range("K13").Value = 17.3731360634624
$ws.Range("L13").Value = 12.95412516594498
$ws.Range("N13").Value = 24.79157597768281

$ws.Range("B14").Value = 22.13057470675116
$ws.Range("C14").Value = 5.31660222422044
$ws.Range("D14").Value = 13.81705824733474
$ws.Range("E14").Value = 13.89154257777728
$ws.Range("G14").Value = 66.22642876803826
$ws.Range("H14").Value = 23.45180359092933
$ws.Range("J14").Value = 8.380848666004081
$ws.Range("K14").Value = 17.35203813309583
$ws.Range("L14").Value = 12.95389566277165
$ws.Range("N14").Value = 24.80035877490794

$ws.Range("B15").Value = 22.11172532543669
$ws.Range("C15").Value = 5.303265615564038
$ws.Range("D15").Value = 13.81460346580405
$ws.Range("E15").Value = 13.89280093920677
$ws.Range("G15").Value = 66.21005212185067
$ws.Range("H15").Value = 23.4533245971013
$ws.Range("J15").Value = 8.38182578858491
$ws.Range("K15").Value = 17.33917488958328
$ws.Range("L15").Value = 12.95379215847522
$ws.Range("N15").Value = 24.8057760449181

$ws.Range("B16").Value = 22.00486267225882
$ws.Range("C16").Value = 5.225926356934089
$ws.Range("D16").Value = 13.80124703470532
$ws.Range("E16").Value = 13.90053534523948
$ws.Range("G16").Value = 66.12001766779967
$ws.Range("H16").Value = 23.46305256511199
$ws.Range("J16").Value = 8.387517229082327
$ws.Range("K16").Value = 17.26639352493245
$ws.Range("L16").Value = 12.95376797533488
$ws.Range("N16").Value = 24.83739607367831

$ws.Range("B17").Value = 21.94036265176475
$ws.Range("C17").Value = 5.177677753887626
$ws.Range("D17").Value = 13.79368988549821
$ws.Range("E17").Value = 13.90574743310722
$ws.Range("G17").Value = 66.06819808486361
$ws.Range("H17").Value = 23.4699231066779
$ws.Range("J17").Value = 8.391090977462452
$ws.Range("K17").Value = 17.22259325587661
$ws.Range("L17").Value = 12.95426159110834
$ws.Range("N17").Value = 24.85730805876794

$ws.Range("B18").Value = 21.90366234371505
$ws.Range("C18").Value = 5.149635600417581
$ws.Range("D18").Value = 13.78957799932096
$ws.Range("E18").Value = 13.90891715663363
$ws.Range("G18").Value = 66.03965137010344
$ws.Range("H18").Value = 23.47420659807681
$ws.Range("J18").Value = 8.393176798795231
$ws.Range("K18").Value = 17.1977188467322
$ws.Range("L18").Value = 12.95473264714067
$ws.Range("N18").Value = 24.86894995892607

$ws.Range("B19").Value = 21.89130561521703
$ws.Range("C19").Value = 5.140091494092871
$ws.Range("D19").Value = 13.78822618657657
$ws.Range("E19").Value = 13.91001990230507
$ws.Range("G19").Value = 66.03020245870253
$ws.Range("H19").Value = 23.47571388764975
$ws.Range("J19").Value = 8.393888232530433
$ws.Range("K19").Value = 17.18935207650011
$ws.Range("L19").Value = 12.95492430624936
$ws.Range("N19").Value = 24.87292418749619

$ws.Range("B20").Value = 21.94718778064583
$ws.Range("C20").Value = 5.182844102880972
$ws.Range("D20").Value = 13.7944700758533
$ws.Range("E20").Value = 13.90517481288498
$ws.Range("G20").Value = 66.07358420087796
$ws.Range("H20").Value = 23.46915739145734
$ws.Range("J20").Value = 8.390707412042016
$ws.Range("K20").Value = 17.22722306501134
$ws.Range("L20").Value = 12.95418968269001
$ws.Range("N20").Value = 24.85516883028039

$ws.Range("B21").Value = 22.13963357320689
$ws.Range("C21").Value = 5.322981200850279
$ws.Range("D21").Value = 13.81824792537739
$ws.Range("E21").Value = 13.89094840959034
$ws.Range("G21").Value = 66.23434911335262
$ws.Range("H21").Value = 23.45109212288263
$ws.Range("J21").Value = 8.380381732978767
$ws.Range("K21").Value = 17.35822267148638
$ws.Range("L21").Value = 12.95395534922205
$ws.Range("N21").Value = 24.79777122903445

$ws.Range("B22").Value = 22.26873991510469
$ws.Range("C22").Value = 5.411887776793231
$ws.Range("D22").Value = 13.83586028293368
$ws.Range("E22").Value = 13.88317792022594
$ws.Range("G22").Value = 66.35053487299649
$ws.Range("H22").Value = 23.44223959832733
$ws.Range("J22").Value = 8.373901494635941
$ws.Range("K22").Value = 17.4465346477816
$ws.Range("L22").Value = 12.95546229498651
$ws.Range("N22").Value = 24.76193967335527

$ws.Range("B23").Value = 22.19953919824316
$ws.Range("C23").Value = 5.364686157618782
$ws.Range("D23").Value = 13.82627159025583
$ws.Range("E23").Value = 13.88718559282155
$ws.Range("G23").Value = 66.28751194846699
$ws.Range("H23").Value = 23.44669407254936
$ws.Range("J23").Value = 8.377335657053321
$ws.Range("K23").Value = 17.39916105020643
$ws.Range("L23").Value = 12.9545064419347
$ws.Range("N23").Value = 24.78090994121421

$ws.Range("B24").Value = 21.94410094877542
$ws.Range("C24").Value = 5.180509340501701
$ws.Range("D24").Value = 13.79411662639662
$ws.Range("E24").Value = 13.9054331548774
$ws.Range("G24").Value = 66.07114525843141
$ws.Range("H24").Value = 23.46950253231774
$ws.Range("J24").Value = 8.39088072468455
$ws.Range("K24").Value = 17.22512897062995
$ws.Range("L24").Value = 12.95422160908817
$ws.Range("N24").Value = 24.8561353703496

$ws.Range("B25").Value = 21.68391669993355
$ws.Range("C25").Value = 4.970792074048349
$ws.Range("D25").Value = 13.76841650297942
$ws.Range("E25").Value = 13.93168836883224
$ws.Range("G25").Value = 65.88590768762538
$ws.Range("H25").Value = 23.50678666273539
$ws.Range("J25").Value = 8.406650536173826
$ws.Range("K25").Value = 17.04965113267352
$ws.Range("L25").Value = 12.96107168930366
$ws.Range("N25").Value = 24.94450782900689
